$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.530.04"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "'1.655.42"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'302.34"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.3838"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").Value = "'0.3600"
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("D9").Value = "'51.10"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").Value = "'0.08201"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("D11").Value = "'1.240"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "'22.39"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "'6.484"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "'7.508"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "'0.00001221"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "'1.654.77"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").Value = "'97.65"
$ws.Range("E18").Value = "  +4.18%  "
$ws.Range("D19").Value = "'0.06985"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "'6.817"
$ws.Range("E20").Value = "  +5.72%  "
$ws.Range("D21").Value = "'17.65"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'12.68"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("D24").Value = "'23.547.02"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'2.513"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "'3.026"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'21.23"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "'152.44"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "'5.230"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'133.99"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'1.841.85"
$ws.Range("E31").Value = "  +3.75%  "
$ws.Range("D32").Value = "'7.216"
$ws.Range("E32").Value = "  +12.12%  "
$ws.Range("D33").Value = "'2.252"
$ws.Range("E33").Value = "  +7.06%  "
$ws.Range("E34").Value = "  +6.08%  "
$ws.Range("D35").Value = "'1.057"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").Value = "'0.02811"
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("D37").Value = "'6.133"
$ws.Range("E37").Value = "  +5.09%  "
$ws.Range("D38").Value = "'0.2499"
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("D39").Value = "'0.08786"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "'0.07017"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("E41").Value = "  +10.65%  "
$ws.Range("D42").Value = "'0.7010"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("D43").Value = "'1.336"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "'15.91"
$ws.Range("E44").Value = "  +5.01%  "
$ws.Range("D45").Value = "'0.6526"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'3.953"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'0.07910"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'128.29"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "'1.192"
$ws.Range("E51").Value = "  +2.22%  "
